# Scheduled-runner update: refresh cached market-board price/profit figures
# (currentAveragePrice*, Leve price & profit columns H:N) across the Leve
# profit sheets, per the latest pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 3704123.8
$ws.Range("I33").Value = 4000411
$ws.Range("K33").Value = 4000411
$ws.Range("M33").Value = -4000182
$ws.Range("H39").Value = 284
$ws.Range("J39").Value = 443.125
$ws.Range("L39").Value = 1329.375
$ws.Range("N39").Value = -1921.375
$ws.Range("H41").Value = 783.61536
$ws.Range("I41").Value = 921.55554
$ws.Range("K41").Value = 921.55554
$ws.Range("M41").Value = -481.55554
$ws.Range("H70").Value = 3495.8
$ws.Range("J70").Value = 5924
$ws.Range("L70").Value = 17772
$ws.Range("N70").Value = -18312
$ws.Range("H73").Value = 3495.8
$ws.Range("J73").Value = 5924
$ws.Range("L73").Value = 17772
$ws.Range("N73").Value = -19644
$ws.Range("H112").Value = 78850.92
$ws.Range("J112").Value = 92884.82000000001
$ws.Range("L112").Value = 278654.46
$ws.Range("N112").Value = -280870.46
$ws.Range("H116").Value = 53831.668
$ws.Range("J116").Value = 10000
$ws.Range("L116").Value = 10000
$ws.Range("N116").Value = -16884
$ws.Range("H127").Value = 1613.6
$ws.Range("J127").Value = 1815.3334
$ws.Range("L127").Value = 5446.0002
$ws.Range("N127").Value = -15366.0002
$ws.Range("H133").Value = 77950
$ws.Range("J133").Value = 77950
$ws.Range("L133").Value = 77950
$ws.Range("N133").Value = -88070
$ws.Range("H137").Value = 1429.5
$ws.Range("I137").Value = 1112.5
$ws.Range("K137").Value = 3337.5
$ws.Range("M137").Value = -787.5
$ws.Range("H138").Value = 3611.7
$ws.Range("J138").Value = 4109.288
$ws.Range("L138").Value = 12327.864
$ws.Range("N138").Value = -22607.864

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4116.4897
$ws.Range("I32").Value = 3604.7556
$ws.Range("K32").Value = 3604.7556
$ws.Range("M32").Value = -3317.7556
$ws.Range("H45").Value = 113425.9
$ws.Range("I45").Value = 125254.445
$ws.Range("K45").Value = 125254.445
$ws.Range("M45").Value = -124877.445
$ws.Range("H61").Value = 2871.814
$ws.Range("J61").Value = 4329.7
$ws.Range("L61").Value = 4329.7
$ws.Range("N61").Value = -4753.7
$ws.Range("H132").Value = 2911.875
$ws.Range("I132").Value = 2670.1072
$ws.Range("K132").Value = 8010.321599999999
$ws.Range("M132").Value = -5480.321599999999
$ws.Range("H133").Value = 98837.89
$ws.Range("J133").Value = 98837.89
$ws.Range("L133").Value = 98837.89
$ws.Range("N133").Value = -103897.89
$ws.Range("H136").Value = 2871.814
$ws.Range("J136").Value = 4329.7
$ws.Range("L136").Value = 12989.1
$ws.Range("N136").Value = -18089.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 74750
$ws.Range("J92").Value = 74750
$ws.Range("L92").Value = 74750
$ws.Range("N92").Value = -79742
$ws.Range("H99").Value = 4696.8125
$ws.Range("I99").Value = 1137.5
$ws.Range("J99").Value = 15374.75
$ws.Range("K99").Value = 1137.5
$ws.Range("L99").Value = 15374.75
$ws.Range("M99").Value = 360.5
$ws.Range("N99").Value = -18370.75
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H134").Value = 2251.9119
$ws.Range("I134").Value = 2259.5454
$ws.Range("K134").Value = 6778.6362
$ws.Range("M134").Value = -4243.6362

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28884.21
$ws.Range("I31").Value = 38977.11
$ws.Range("J31").Value = 4110.727
$ws.Range("K31").Value = 38977.11
$ws.Range("L31").Value = 4110.727
$ws.Range("M31").Value = -38682.11
$ws.Range("N31").Value = -4700.727
$ws.Range("H34").Value = 28884.21
$ws.Range("I34").Value = 38977.11
$ws.Range("J34").Value = 4110.727
$ws.Range("K34").Value = 38977.11
$ws.Range("L34").Value = 4110.727
$ws.Range("M34").Value = -38775.11
$ws.Range("N34").Value = -4514.727
$ws.Range("H44").Value = 47500
$ws.Range("J44").Value = 80000
$ws.Range("L44").Value = 80000
$ws.Range("N44").Value = -80884
$ws.Range("H99").Value = 7629.75
$ws.Range("I99").Value = 2752.5
$ws.Range("J99").Value = 12507
$ws.Range("K99").Value = 2752.5
$ws.Range("L99").Value = 12507
$ws.Range("M99").Value = -1254.5
$ws.Range("N99").Value = -15503
$ws.Range("H126").Value = 7629.75
$ws.Range("I126").Value = 2752.5
$ws.Range("J126").Value = 12507
$ws.Range("K126").Value = 8257.5
$ws.Range("L126").Value = 37521
$ws.Range("M126").Value = -5787.5
$ws.Range("N126").Value = -42461
$ws.Range("H132").Value = 8066.625
$ws.Range("I132").Value = 5092.625
$ws.Range("J132").Value = 14014.625
$ws.Range("K132").Value = 15277.875
$ws.Range("L132").Value = 42043.875
$ws.Range("M132").Value = -12747.875
$ws.Range("N132").Value = -47103.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 72.44444
$ws.Range("I2").Value = 70.59999999999999
$ws.Range("K2").Value = 423.6
$ws.Range("M2").Value = -310.6
$ws.Range("H23").Value = 1357.5714
$ws.Range("I23").Value = 2166.3333
$ws.Range("K23").Value = 6498.999899999999
$ws.Range("M23").Value = -6263.999899999999
$ws.Range("H114").Value = 730.2308
$ws.Range("I114").Value = 759.2222
$ws.Range("J114").Value = 665
$ws.Range("K114").Value = 2277.6666
$ws.Range("L114").Value = 1995
$ws.Range("M114").Value = 976.3334
$ws.Range("N114").Value = -8503
$ws.Range("H117").Value = 1264.75
$ws.Range("J117").Value = 3000
$ws.Range("L117").Value = 9000
$ws.Range("N117").Value = -15884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 794.94446
$ws.Range("J23").Value = 947.2308
$ws.Range("L23").Value = 947.2308
$ws.Range("N23").Value = -1393.2308
$ws.Range("H113").Value = 2015.1875
$ws.Range("I113").Value = 2018.7693
$ws.Range("J113").Value = 1999.6666
$ws.Range("K113").Value = 2018.7693
$ws.Range("L113").Value = 1999.6666
$ws.Range("M113").Value = 151.2307000000001
$ws.Range("N113").Value = -6339.6666
$ws.Range("H126").Value = 3624.8057
$ws.Range("I126").Value = 3187.4167
$ws.Range("K126").Value = 9562.250100000001
$ws.Range("M126").Value = -7092.250100000001
$ws.Range("H132").Value = 3998.0476
$ws.Range("I132").Value = 3831.0557
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 11493.1671
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -8963.167099999999
$ws.Range("N132").Value = -20060
$ws.Range("H141").Value = 81394.39999999999
$ws.Range("J141").Value = 81394.39999999999
$ws.Range("L141").Value = 81394.39999999999
$ws.Range("N141").Value = -91754.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 15860551
$ws.Range("J23").Value = 50000000
$ws.Range("L23").Value = 50000000
$ws.Range("N23").Value = -50000460
$ws.Range("H40").Value = 6245.125
$ws.Range("I40").Value = 5994.4287
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 5994.4287
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -5858.4287
$ws.Range("N40").Value = -8272
$ws.Range("H132").Value = 4815.231
$ws.Range("I132").Value = 4959.2
$ws.Range("J132").Value = 4335.3335
$ws.Range("K132").Value = 14877.6
$ws.Range("L132").Value = 13006.0005
$ws.Range("M132").Value = -12347.6
$ws.Range("N132").Value = -18066.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 20108.666
$ws.Range("J74").Value = 20108.666
$ws.Range("L74").Value = 20108.666
$ws.Range("N74").Value = -21980.666
$ws.Range("H77").Value = 20108.666
$ws.Range("J77").Value = 20108.666
$ws.Range("L77").Value = 60325.99800000001
$ws.Range("N77").Value = -69685.99800000001
$ws.Range("H100").Value = 1393.421
$ws.Range("J100").Value = 1613.1428
$ws.Range("L100").Value = 3226.2856
$ws.Range("N100").Value = -4308.2856
$ws.Range("H122").Value = 2759.75
$ws.Range("I122").Value = 2642.2188
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 7926.6564
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -5476.6564
$ws.Range("N122").Value = -16000
$ws.Range("H132").Value = 2412.6365
$ws.Range("I132").Value = 2380
$ws.Range("J132").Value = 2499.6667
$ws.Range("K132").Value = 7140
$ws.Range("L132").Value = 7499.000100000001
$ws.Range("M132").Value = -4610
$ws.Range("N132").Value = -12559.0001
$ws.Range("H137").Value = 100694.5
$ws.Range("J137").Value = 100694.5
$ws.Range("L137").Value = 100694.5
$ws.Range("N137").Value = -110894.5
